$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 53
$ws1.Range("F5").Value = 5019
$ws1.Range("F7").Value = 86
$ws1.Range("F9").Value = 50

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 53
$ws4.Range("F9").Value = 5019
$ws4.Range("F11").Value = 86
$ws4.Range("F14").Value = 50
